$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: sending cluster FAPs, ligand Agt, receptor Agtr2, target cluster FAPs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Agt"
$ws.Range("C2").Value = "Agtr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4683593333333333
$ws.Range("H2").Value = 1.405078
$ws.Range("I2").Value = 0.6051469521021553
$ws.Range("J2").Value = 0.6051469521021552
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 64.154275
$ws.Range("N2").Value = 192.462825
$ws.Range("O2").Value = 0.9711753522845754
$ws.Range("P2").Value = 0.9711753522845754
$ws.Range("Q2").Value = 30.04725346948333
$ws.Range("R2").Value = 270.42528122535
$ws.Range("S2").Value = 0.5877038043917477
$ws.Range("T2").Value = 0.5877038043917476

# Row 3: sending cluster FAPs, ligand Agt, receptor Agtr2, target cluster sCs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Agt"
$ws.Range("C3").Value = "Agtr2"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.4683593333333333
$ws.Range("H3").Value = 1.405078
$ws.Range("I3").Value = 0.6051469521021553
$ws.Range("J3").Value = 0.6051469521021552
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.904109666666667
$ws.Range("N3").Value = 5.712329
$ws.Range("O3").Value = 0.02882464771542451
$ws.Range("P3").Value = 0.02882464771542451
$ws.Range("Q3").Value = 0.8918075340735556
$ws.Range("R3").Value = 8.026267806662
$ws.Range("S3").Value = 0.0174431477104075
$ws.Range("T3").Value = 0.01744314771040749

# Row 4: sending cluster sCs, ligand Agt, receptor Agtr2, target cluster FAPs
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Agt"
$ws.Range("C4").Value = "Agtr2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3056003333333333
$ws.Range("H4").Value = 0.916801
$ws.Range("I4").Value = 0.3948530478978448
$ws.Range("J4").Value = 0.3948530478978448
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 64.154275
$ws.Range("N4").Value = 192.462825
$ws.Range("O4").Value = 0.9711753522845754
$ws.Range("P4").Value = 0.9711753522845754
$ws.Range("Q4").Value = 19.60556782475833
$ws.Range("R4").Value = 176.450110422825
$ws.Range("S4").Value = 0.3834715478928278
$ws.Range("T4").Value = 0.3834715478928278

# Row 5: sending cluster sCs, ligand Agt, receptor Agtr2, target cluster sCs
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Agt"
$ws.Range("C5").Value = "Agtr2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.3056003333333333
$ws.Range("H5").Value = 0.916801
$ws.Range("I5").Value = 0.3948530478978448
$ws.Range("J5").Value = 0.3948530478978448
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.904109666666667
$ws.Range("N5").Value = 5.712329
$ws.Range("O5").Value = 0.02882464771542451
$ws.Range("P5").Value = 0.02882464771542451
$ws.Range("Q5").Value = 0.5818965488365555
$ws.Range("R5").Value = 5.237068939529
$ws.Range("S5").Value = 0.01138150000501702
$ws.Range("T5").Value = 0.01138150000501702
